$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row is row 73, one past the previous last data row (72).
$newRow = 73

# Replicate the cell formatting (styles) used on the row above for the
# two columns that carry an explicit style: A (bold/bordered index
# column) and E (date-time column).
$ws.Range("A72").Copy()
$ws.Range("A73").PasteSpecial(-4122)
$ws.Range("E72").Copy()
$ws.Range("E73").PasteSpecial(-4122)

# Populate the new row's values.
$ws.Range("A73").Value = 72
$ws.Range("B73").Value = "bosnia-and-herzegovina"
$ws.Range("C73").Value = "premijer-liga-bih"
$ws.Range("D73").Value = "2023-2024"
$ws.Range("E73").Value = 45233.75
$ws.Range("F73").Value = "GOSK Gabela"
$ws.Range("G73").Value = 0
$ws.Range("H73").Value = "Velez Mostar"
$ws.Range("I73").Value = 4
$ws.Range("J73").Value = 2.89
$ws.Range("K73").Value = "02/11/2023 06:12"
$ws.Range("L73").Value = 3.27
$ws.Range("M73").Value = "03/11/2023 17:56"
$ws.Range("N73").Value = 3.04
$ws.Range("O73").Value = "02/11/2023 06:12"
$ws.Range("P73").Value = 3.28
$ws.Range("Q73").Value = "03/11/2023 17:56"
$ws.Range("R73").Value = 2.27
$ws.Range("S73").Value = "02/11/2023 06:12"
$ws.Range("T73").Value = 2.17
$ws.Range("U73").Value = "03/11/2023 17:56"
$ws.Range("V73").Value = "https://www.betexplorer.com/football/bosnia-and-herzegovina/premijer-liga-bih/nk-gosk-gabela-velez-mostar/je0Q5Eg5/"
